# Update generator linear examples: refresh the randomly generated
# expression/value data on the "Restricciones_del_follower",
# "Punto_modificado", "Vector_bf", "Vector_BF" and "Vector_Alpha" sheets.
#
# Most of the written values are decimal numbers that are stored as TEXT
# in the workbook (they live in columns together with algebraic
# expressions such as "-3 - x + 1.9035414930053012y"). Setting a range's
# NumberFormat to "@" (Text) before assigning the value keeps Excel from
# auto-converting the numeric-looking strings to real numbers; clearing
# the format again afterwards drops the leftover style so the cell is
# left exactly as it was (General / default style) but still holding a
# text value.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# --- Restricciones_del_follower -----------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$ws.Range("A2").Value = "1.2243551978833835y"
Set-TextValue $ws.Range("B2") "5.356092504467225"
$ws.Range("C2").Value = "J_0_L0_v"
Set-TextValue $ws.Range("D2") "0.33468162538227564"
Set-TextValue $ws.Range("E2") "0.7610371704678504"
Set-TextValue $ws.Range("F2") "0.5409063754085729"

$ws.Range("A3").Value = "-3 - x + 1.9035414930053012y"
Set-TextValue $ws.Range("B3") "1.142384129045567"
$ws.Range("C3").Value = "J_0_L0_v"
Set-TextValue $ws.Range("D3") "0.9092567913461869"
Set-TextValue $ws.Range("E3") "0.0052664810917755125"
Set-TextValue $ws.Range("F3") "0.8409632524951235"

$ws.Range("A4").Value = "-12 + x + 1.1164224624624124y"
Set-TextValue $ws.Range("B4") "-2.9311801143896368"
Set-TextValue $ws.Range("D4") "0.7906785535517057"
Set-TextValue $ws.Range("E4") "0.1970156715059802"
Set-TextValue $ws.Range("F4") "0.629307423434291"

$ws.Range("A5").Value = "-12 + 4x + 1.924870762754641y"
Set-TextValue $ws.Range("B5") "13.16015372672777"
$ws.Range("C5").Value = "J_Ne_L0_v"
Set-TextValue $ws.Range("D5") "0.5618257705012442"
Set-TextValue $ws.Range("E5") "0.1563052103601904"
Set-TextValue $ws.Range("F5") "0.8503862843164217"

# --- Punto_modificado -----------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "4.184892416399492"
Set-TextValue $ws.Range("B2") "4.374623078112156"

# --- Vector_bf --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws.Range("A2") "-5.1047505149695365"

# --- Vector_BF --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $ws.Range("A2") "0.1830299681450337"
Set-TextValue $ws.Range("A3") "1.537375168659966"

# --- Vector_Alpha ---------------------------------------------------------
# This one is a genuine numeric cell (no text type in the original file).
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.4417887687687938
